# tcpip.docx edits per commit diff

$d = $word.ActiveDocument

# 1. "protocol de- signs" -> "protocol designs"
$d.Content.Find.Execute("protocol de- signs", $true, $false, $false, $false, $false,
                         $true, 1, $false, "protocol designs", 2) | Out-Null

# 2. "# sysctl -q net.ipv4.tcp_max_syn_backlog" -> "sudo sysctl -q net.ipv4.tcp_max_syn_backlog"
$d.Content.Find.Execute("# sysctl -q net.ipv4.tcp_max_syn_backlog", $true, $false, $false, $false, $false,
                         $true, 1, $false, "sudo sysctl -q net.ipv4.tcp_max_syn_backlog", 2) | Out-Null

# 3. typo fix + insertion: "to send one packaet at a time, and" -> "to send one packet at a time to that IP address, and"
$d.Content.Find.Execute([string]([char]0x201C) + "-c 1" + [string]([char]0x201D) + " to send one packaet at a time, and",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         [string]([char]0x201C) + "-c 1" + [string]([char]0x201D) + " to send one packet at a time to that IP address, and",
                         2) | Out-Null

# 4. "whether the attack is successful or not.  " -> "whether the attack has potential to succeed."
$d.Content.Find.Execute("whether the attack is successful or not.  ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "whether the attack has potential to succeed.", 2) | Out-Null

# 5. "# sysctl -a | grep cookie" -> "sudo sysctl -a | grep cookie"
$d.Content.Find.Execute("# sysctl -a | grep cookie", $true, $false, $false, $false, $false,
                         $true, 1, $false, "sudo sysctl -a | grep cookie", 2) | Out-Null

# 6. "# sysctl -w net.ipv4.tcp_syncookies=0" -> "sudo sysctl -w net.ipv4.tcp_syncookies=0"
$d.Content.Find.Execute("# sysctl -w net.ipv4.tcp_syncookies=0", $true, $false, $false, $false, $false,
                         $true, 1, $false, "sudo sysctl -w net.ipv4.tcp_syncookies=0", 2) | Out-Null

# 7. "# sysctl -w net.ipv4.tcp_syncookies=1" -> "sudo sysctl -w net.ipv4.tcp_syncookies=1"
$d.Content.Find.Execute("# sysctl -w net.ipv4.tcp_syncookies=1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "sudo sysctl -w net.ipv4.tcp_syncookies=1", 2) | Out-Null

# 8. "Relative Sequence Number and Window Scaling" -> "Relative Sequence Number"
$d.Content.Find.Execute("Relative Sequence Number and Window Scaling", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Relative Sequence Number", 2) | Out-Null

# 9. "3.4 Task 4 : TCP Session Hijacking" -> "3.3 Task 3 : TCP Session Hijacking"
$d.Content.Find.Execute("3.4 Task 4 : TCP Session Hijacking", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.3 Task 3 : TCP Session Hijacking", 2) | Out-Null

# 10. "npig" -> "nping"
$d.Content.Find.Execute("npig", $true, $false, $false, $false, $false,
                         $true, 1, $false, "nping", 2) | Out-Null

# 11. "3.5 Task 5 : Creating Reverse Shell using TCP Session Hijacking" -> "3.4 Task 4 : Creating Reverse Shell using TCP Session Hijacking"
$d.Content.Find.Execute("3.5 Task 5 : Creating Reverse Shell using TCP Session Hijacking", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.4 Task 4 : Creating Reverse Shell using TCP Session Hijacking", 2) | Out-Null

# 12. insert three new paragraphs after the "...has potential to succeed." paragraph
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*has potential to succeed.*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()

    $p1 = $d.Paragraphs($i + 1)
    $p1.Range.Text = "To make your attack easier to succeed, we will shrink the size of the backlog queue to 5:"

    $p1.Range.InsertParagraphAfter()
    $p2 = $d.Paragraphs($i + 2)
    $p2.Range.Text = "sudo sysctl -w net.ipv4.tcp_max_syn_backlog=5"
    $p2.Format.FirstLineIndent = 17.05

    $p2.Range.InsertParagraphAfter()
    $p3 = $d.Paragraphs($i + 3)
    $p3.Range.Text = "Send five packets via nping and then try to telnet to the server via the user component.  Report on your success."
    $p3.Format.FirstLineIndent = -0.01
}

Write-Output "phase1 ok"
